# Add a "Serial" column (O) to the calibration data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("O1").Value = "Serial"

# First data value is a literal 1
$ws.Range("O2").Value = 1

# Row 3 uses a plain (non-shared) formula referencing the row above
$ws.Range("O3").Formula = "=1+O2"

# Rows 4:10 use a shared formula referencing the row above
$ws.Range("O4:O10").Formula = "=1+O3"

# Update the selected cell to match the target state
$ws.Range("O5").Select()
